# "Added outsourcing to simulation"
# Duplicates the Resource Staffing data/chart section (rows 1-6, charts 1 & 2)
# into a new section at rows 50-55 with a new heading and new (smaller)
# staffing numbers, and adds two new charts mirroring the first pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block (rows 50-55), mirrors rows 1-6 with new values ---
$ws.Range("A50").Value = "Updated with New Resources"

$ws.Range("A51").Value = "With Error"
$ws.Range("C51").Value = "Without Error"

$ws.Range("A52").Value = "Designers"
$ws.Range("B52").Value = 8
$ws.Range("C52").Value = "Designers"
$ws.Range("D52").Value = 8

$ws.Range("A53").Value = "Implementers"
$ws.Range("B53").Value = 8
$ws.Range("C53").Value = "Implementers"
$ws.Range("D53").Value = 8

$ws.Range("A54").Value = "Testers"
$ws.Range("B54").Value = 8
$ws.Range("C54").Value = "Testers"
$ws.Range("D54").Value = 8

$ws.Range("A55").Value = "Managers"
$ws.Range("B55").Value = 1
$ws.Range("C55").Value = "Managers"
$ws.Range("D55").Value = 1

# --- Column A widened to fit the new, longer label ---
$ws.Columns.Item(1).ColumnWidth = 24.5

# --- New chart #1 ("w/ Error"), mirrors Chart 2 (chart1.xml) but on A52:B55 ---
$co3 = $ws.ChartObjects().Add(61, 873, 369.4794921875, 257)
$co3.Name = "Chart 4"
$chart3 = $co3.Chart
$chart3.ChartType = 52
$chart3.SetSourceData($ws.Range("A52:B55"))
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Resource Staffing Chart (w/ Error)"
$chart3.HasLegend = $false

$catAx3 = $chart3.Axes(1)
$catAx3.HasTitle = $true
$catAx3.AxisTitle.Text = "Employee Type"

$valAx3 = $chart3.Axes(2)
$valAx3.HasTitle = $true
$valAx3.AxisTitle.Text = "# of People"
$valAx3.HasMajorGridlines = $true

$ser3 = $chart3.SeriesCollection(1)
$ser3.HasDataLabels = $false

$chart3.ChartGroups(1).GapWidth = 150
$chart3.ChartGroups(1).Overlap = 100

# --- New chart #2 ("w/o Error"), mirrors Chart 3 (chart2.xml) but on C52:D55 ---
$co4 = $ws.ChartObjects().Add(61, 1148, 365.4794921875, 266)
$co4.Name = "Chart 5"
$chart4 = $co4.Chart
$chart4.ChartType = 52
$chart4.SetSourceData($ws.Range("C52:D55"))
$chart4.HasTitle = $true
$chart4.ChartTitle.Text = "Resource Staffing Chart (w/o Error)"
$chart4.HasLegend = $false

$catAx4 = $chart4.Axes(1)
$catAx4.HasTitle = $true
$catAx4.AxisTitle.Text = "Employee Type"

$valAx4 = $chart4.Axes(2)
$valAx4.HasTitle = $true
$valAx4.AxisTitle.Text = "# of People"
$valAx4.HasMajorGridlines = $true

$ser4 = $chart4.SeriesCollection(1)
$ser4.HasDataLabels = $false

$chart4.ChartGroups(1).GapWidth = 150
$chart4.ChartGroups(1).Overlap = 100

# --- Scroll/selection bookkeeping to match where the user ended up ---
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$null = $ws.Range("H89").Select()
